$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the trigger/download values in column C (rows 3, 7, 10)
$ws.Range("C3").Value = 5
$ws.Range("C7").Value = 2
$ws.Range("C10").Value = 2

# Move/record the active selection at C11
$ws.Range("C11").Select()
